# Move raw data in folder
# This workbook (Result/Data_Fig3.xlsx) had its single sheet renamed from
# "Fig4" to "Fig3" to match the file name / figure numbering after the raw
# data files were reorganised into the Result folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet to match the figure it now represents ---
$ws.Name = "Fig3"

# --- Drop the (visually inert) placeholder style that a handful of cells
# in the small regression-summary block (columns R:S) were still carrying.
# That style has no fill/border/number-format of its own (fontId=0,
# fillId=0, borderId=0) so clearing direct formatting on those cells is a
# no-op visually but removes the redundant style reference, same as the
# cleanup that happened when the workbook was last resaved.
$ws.Range("R5:S8").ClearFormats()
$ws.Range("R15:S18").ClearFormats()
$ws.Range("R25:S28").ClearFormats()
